# Atualização automática: 2025-08-27 21:00:26
# Corrects the First_Detection_Image filename and the bounding-box /
# confidence values for rows 16 and 17 on the active sheet.
#
# Columns: D = First_Detection_Image, I = First_Coords, J = First_Confidence
# I/J values look numeric (comma-delimited numbers / decimals), so a
# leading apostrophe is used to force them to stay plain text, exactly
# like typing them into Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I16").Value = "'642,530,686,576"

# Row 17
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I17").Value = "'794,481,831,526"
$ws.Range("J17").Value = "'0.71"
